$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - 裂变弹匣 (ballistic tier 1)
$ws.Range("H6").Value = "multiShot:2|damage:+6"
$ws.Range("J6").Value = "触发时连续射出两轮弹幕。"
$ws.Range("K6").Value = "icons/skill/ballistic-tier1.svg"

# Row 7 - 轨迹稳流 (ballistic tier 2)
$ws.Range("H7").Value = "multiShot:3|stability:+12|projectileSize:+12|multiShotAngle:4"
$ws.Range("J7").Value = "连续三轮射击并稳定弹道。"
$ws.Range("K7").Value = "icons/skill/ballistic-tier2.svg"

# Row 8 - 深域贯穿 (ballistic tier 3)
$ws.Range("H8").Value = "multiShot:4|pierce:+1|ricochet:+1"
$ws.Range("J8").Value = "连续四轮射击并令子弹贯穿弹射。"
$ws.Range("K8").Value = "icons/skill/ballistic-tier3.svg"

# Row 9 - 棱镜导光 (energy tier 1) - icon only
$ws.Range("K9").Value = "icons/skill/energy-tier1.svg"

# Row 10 - 谐振折叠 (energy tier 2) - icon only
$ws.Range("K10").Value = "icons/skill/energy-tier2.svg"

# Row 11 - 相干放射 (energy tier 3) - icon only
$ws.Range("K11").Value = "icons/skill/energy-tier3.svg"

# Row 12 - 相位壁垒 (guardian tier 1) - icon only
$ws.Range("K12").Value = "icons/skill/guardian-tier1.svg"

# Row 13 - 护盾崩击 (guardian tier 2) - icon only
$ws.Range("K13").Value = "icons/skill/guardian-tier2.svg"

# Row 14 - 寂光回响 (guardian tier 3) - icon only
$ws.Range("K14").Value = "icons/skill/guardian-tier3.svg"

# Row 15 - 术式镀层 (workshop tier 1)
$ws.Range("H15").Value = "projectileSize:+24|elementSlow:+18|elementSlowDuration:+1.2"
$ws.Range("K15").Value = "icons/skill/workshop-tier1.svg"

# Row 16 - 弹道精铸 (workshop tier 2) - icon only
$ws.Range("K16").Value = "icons/skill/workshop-tier2.svg"

# Row 17 - 秘火迸流 (workshop tier 3) - icon only
$ws.Range("K17").Value = "icons/skill/workshop-tier3.svg"
